# Fixed TWOI to SWOI translation
# Missing parantheses, time window was off by about 8 ms, redid tables & graphs
#
# This script updates the statistical values in the LPP exploratory table
# to reflect the corrected time-window analysis.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Row: across groups: CS+av vs CS+neu
Replace-Exact "3.09" "3.11"
Replace-Exact "19.729" "20.748"

# Row: across groups: CS+av vs CS-
Replace-Exact "3.23" "3.25"
Replace-Exact "27.988" "29.108"

# Row: across groups: CSneu vs CS-
Replace-Exact ".919" ".923"
Replace-Exact " 0.158" " 0.157"

# Row: imagery: CS+av vs CS+neu
Replace-Exact "1.68" "1.69"
Replace-Exact ".053" ".052"
Replace-Exact "0.34" "0.35"
Replace-Exact " 1.373" " 1.396"

# Row: imagery: CS+av vs CS-
Replace-Exact "3.59" "3.60"
Replace-Exact "0.73" "0.74"
Replace-Exact "47.347" "49.322"

# Row: imagery: CSneu vs CS-
Replace-Exact "1.44" "1.47"
Replace-Exact ".163" ".156"
Replace-Exact "0.29" "0.30"
Replace-Exact " 0.534" " 0.551"

# Row: classical: CS+av vs CS+neu
Replace-Exact "2.61" "2.63"
Replace-Exact ".008" ".007"
Replace-Exact "0.53" "0.54"
Replace-Exact " 6.606" " 6.843"

# Row: classical: CS+av vs CS-
Replace-Exact " 1.092" " 1.096"

# Row: classical: CSneu vs CS-
Replace-Exact "-1.33" "-1.35"
Replace-Exact ".198" ".191"
Replace-Exact "-0.27" "-0.28"
Replace-Exact " 0.466" " 0.478"

# Row: between groups: delta CS+av / CS+neu
Replace-Exact "1.18" "1.20"
Replace-Exact ".734" ".709"
Replace-Exact " 0.506" " 0.516"

# Row: between groups: delta CS+av / CS-
Replace-Exact "-0.64" "-0.65"
Replace-Exact "-0.18" "-0.19"
Replace-Exact " 0.339" " 0.342"

# Row: between groups: delta CSneu / CS-
Replace-Exact "-1.96" "-1.99"
Replace-Exact ".169" ".158"
Replace-Exact "-0.56" "-0.57"
Replace-Exact " 1.330" " 1.400"
